# Apply updates to "Principles of Computational Modelling in Neuroscience.xlsx"
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# --- Fill in the new evaluation rows (7-10) on Sheet1 ---
$rowsData = @(
    @{ Row = 7;  A = 42965.033333333333; B = 267; C = 42965.073611111111; D = 277 },
    @{ Row = 8;  A = 42965.666666666664; B = 278; C = 42965.681944444441; D = 285 },
    @{ Row = 9;  A = 42965.931944444441; B = 286; C = 42965.956944444442; D = 296 },
    @{ Row = 10; A = 42965.821527777778; B = 297; C = 42965.859722222223; D = 313 }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $sheet1.Cells.Item($row, 1).Value = $r.A
    $sheet1.Cells.Item($row, 2).Value = $r.B
    $sheet1.Cells.Item($row, 3).Value = $r.C
    $sheet1.Cells.Item($row, 4).Value = $r.D
    $sheet1.Cells.Item($row, 5).Formula = "=C$row-A$row"
    $sheet1.Cells.Item($row, 6).Formula = "=D$row-B$row+1"
    $sheet1.Cells.Item($row, 7).Formula = "=F$row/(E$row*24*60)"
    $sheet1.Cells.Item($row, 8).Formula = "=G$row*60"
    $sheet1.Cells.Item($row, 9).Formula = "=Sheet2!`$B`$25/MEDIAN(`$H`$2:H$row)"
    $sheet1.Cells.Item($row, 10).Formula = "=Sheet2!`$B`$25/AVERAGE(`$H`$2:H$row)"

    # Columns B and D carried an inherited "Normal w/ explicit font" style from
    # the blank placeholder rows; clear it back to the plain default like the
    # other data rows (2-9) already have.
    $sheet1.Cells.Item($row, 2).Style = "Normal"
    $sheet1.Cells.Item($row, 4).Style = "Normal"
}

# --- Recalculate so cached formula results refresh (Sheet2 stats, chart caches) ---
$excel.CalculateFullRebuild()

# Best-effort: nudge the line chart on Sheet2 to re-pull its cached data
# points now that Sheet1!I7:J10 are populated (harmless if unsupported).
try {
    $co = $sheet2.ChartObjects().Item(1)
    $co.Chart.Refresh()
} catch {
}

# --- Fix up sheet selection / active tab state ---
# Sheet1 becomes the active sheet (was Sheet2), with C11 selected.
# Sheet2 keeps its existing B26 selection untouched; it just stops being the
# active tab once Sheet1 is activated.
$sheet1.Activate()
$sheet1.Range("C11").Select()
